$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Hdfc" / "Hdfc Bank" entries with "Canara" / "Canara Bank"
$ws.Range("A2").Value = "Canara"
$ws.Range("B2").Value = "Canara Bank"

# Update selection to match the new active cell
$ws.Range("B2").Select()
